$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.07292628288269
$ws.Range("B1").Value = 1.293322324752808
$ws.Range("C1").Value = 1.735774993896484
$ws.Range("D1").Value = 3.264542818069458
$ws.Range("E1").Value = 2.308448314666748
